# #5: property boat&car done
# Extends the "汽車" (car) sheet (3rd worksheet) with the common trailing
# columns (property_category .. index) used by the other property sheets,
# and turns row 1 into the shared header row used across the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Grow the header row (row 1) with the extra header cells, copying the
#     style from the existing header cell B1 so the new cells pick up the
#     same formatting (style index) as the rest of the header row.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("B1").Copy($ws.Range("G1"))
$ws.Range("B1").Copy($ws.Range("H1"))
$ws.Range("B1").Copy($ws.Range("I1"))
$ws.Range("B1").Copy($ws.Range("J1"))
$ws.Range("B1").Copy($ws.Range("K1"))
$ws.Range("B1").Copy($ws.Range("L1"))
$ws.Range("B1").Copy($ws.Range("M1"))
$ws.Range("B1").Copy($ws.Range("N1"))

# Row 1 becomes the common header row shared by every property sheet.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Grow the data row (row 2) with the extra trailing data cells, copying
#     the style from the existing data cell B2.
$ws.Range("B2").Copy($ws.Range("H2"))
$ws.Range("B2").Copy($ws.Range("I2"))
$ws.Range("B2").Copy($ws.Range("J2"))
$ws.Range("B2").Copy($ws.Range("K2"))
$ws.Range("B2").Copy($ws.Range("L2"))
$ws.Range("B2").Copy($ws.Range("M2"))
$ws.Range("B2").Copy($ws.Range("N2"))

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("K2").Value = "鄭麗君"
$ws.Range("L2").Value = 1764
$ws.Range("M2").Value = "tmp81521"
$ws.Range("N2").Value = 49
